$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: header + value
$ws.Range("D1").Value = "Fjerde kolonne"
$ws.Range("D2").Value = "Der bringer ændringer i sin egen branch"

# Make column D wide enough (Excel drags column width to the pixel
# then reports it back in character units)
$ws.Columns.Item(4).ColumnWidth = 40.5

# New small column F, repeating the first header ("Data") and a list of numbers
$ws.Range("F1").Value = "Data"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 4

# Move the active selection the way the author left it
$ws.Range("F6").Select()

$wb.Windows.Item(1).Left = 3160
$wb.Windows.Item(1).Top = 3760
